$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 275 (Feria Lagunitas de
# Puerto Montt - Perejil), pushing the existing rows 275-346 down to 276-347.
$ws.Rows.Item(275).Insert()

$ws.Cells.Item(275, 1).Value = 4
$ws.Cells.Item(275, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(275, 3).Value = "Los Lagos"
$ws.Cells.Item(275, 4).Value = 44964
$ws.Cells.Item(275, 5).Value = 10
$ws.Cells.Item(275, 6).Value = 100112044
$ws.Cells.Item(275, 7).Value = "Perejil"
$ws.Cells.Item(275, 8).Value = "Sin especificar"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 120
$ws.Cells.Item(275, 11).Value = 7000
$ws.Cells.Item(275, 12).Value = 7000
$ws.Cells.Item(275, 13).Value = 7000
$ws.Cells.Item(275, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(275, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(275, 16).Value = 3500
$ws.Cells.Item(275, 17).Value = 2
$ws.Cells.Item(275, 18).Value = "Hortaliza"
